$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date for each row (rows 2-151).
# Update it from 45171 (2023-09-02) to 45172 (2023-09-03) for every data row.
$ws.Range("C2:C151").Value = 45172
